$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.747.37"
Set-TextValue $ws.Range("E2") "  +0.58%  "
Set-TextValue $ws.Range("D3") "1.594.85"
Set-TextValue $ws.Range("E3") "  -0.21%  "
Set-TextValue $ws.Range("E4") "  +0.14%  "
Set-TextValue $ws.Range("D5") "209.67"
Set-TextValue $ws.Range("E5") "  +0.35%  "
Set-TextValue $ws.Range("D6") "0.502"
Set-TextValue $ws.Range("E6") "  +0.28%  "
Set-TextValue $ws.Range("E7") "  +0.09%  "
Set-TextValue $ws.Range("E8") "  -0.04%  "
Set-TextValue $ws.Range("D9") "0.254"
Set-TextValue $ws.Range("E9") "  +0.69%  "
Set-TextValue $ws.Range("E10") "  +0.37%  "
Set-TextValue $ws.Range("D11") "0.0869"
Set-TextValue $ws.Range("E11") "  -0.39%  "
Set-TextValue $ws.Range("D12") "1.821.19"
Set-TextValue $ws.Range("E12") "  -0.20%  "
Set-TextValue $ws.Range("D13") "1.610.36"
Set-TextValue $ws.Range("E13") "  +0.51%  "
Set-TextValue $ws.Range("E14") "  -0.47%  "
Set-TextValue $ws.Range("D15") "0.532"
Set-TextValue $ws.Range("E15") "  -1.79%  "
Set-TextValue $ws.Range("D16") "27.741.57"
Set-TextValue $ws.Range("E16") "  +0.59%  "
Set-TextValue $ws.Range("D17") "63.38"
Set-TextValue $ws.Range("E17") "  -0.45%  "
Set-TextValue $ws.Range("D18") "218.96"
Set-TextValue $ws.Range("E18") "  +0.39%  "
Set-TextValue $ws.Range("E19") "  +0.92%  "
Set-TextValue $ws.Range("E20") "  -1.00%  "
Set-TextValue $ws.Range("E21") "  +0.18%  "
Set-TextValue $ws.Range("D22") "4.16"
Set-TextValue $ws.Range("E22") "  -0.99%  "
Set-TextValue $ws.Range("D23") "9.80"
Set-TextValue $ws.Range("E23") "  +0.44%  "
Set-TextValue $ws.Range("D24") "1.98"
Set-TextValue $ws.Range("E24") "  -2.02%  "
Set-TextValue $ws.Range("D25") "153.86"
Set-TextValue $ws.Range("E25") "  -0.33%  "
Set-TextValue $ws.Range("D26") "7.12"
Set-TextValue $ws.Range("E26") "  +5.89%  "
Set-TextValue $ws.Range("E27") "  +0.15%  "
Set-TextValue $ws.Range("D28") "15.16"
Set-TextValue $ws.Range("E28") "  +0.83%  "
Set-TextValue $ws.Range("E29") "  -0.03%  "
Set-TextValue $ws.Range("E30") "  +0.43%  "
Set-TextValue $ws.Range("D31") "0.0477"
Set-TextValue $ws.Range("E31") "  +2.14%  "
Set-TextValue $ws.Range("D32") "3.23"
Set-TextValue $ws.Range("E32") "  -2.22%  "
Set-TextValue $ws.Range("D33") "1.383.78"
Set-TextValue $ws.Range("E33") "  +0.94%  "
Set-TextValue $ws.Range("E34") "  +0.83%  "
Set-TextValue $ws.Range("E35") "  -0.16%  "
Set-TextValue $ws.Range("D36") "0.969"
Set-TextValue $ws.Range("E36") "  +0.94%  "
Set-TextValue $ws.Range("E37") "  +0.83%  "
Set-TextValue $ws.Range("E38") "  +3.02%  "
Set-TextValue $ws.Range("E39") "  +0.17%  "
Set-TextValue $ws.Range("D40") "0.830"
Set-TextValue $ws.Range("E40") "  +1.73%  "
Set-TextValue $ws.Range("E41") "  +0.10%  "
Set-TextValue $ws.Range("E42") "  +0.53%  "
Set-TextValue $ws.Range("D43") "64.64"
Set-TextValue $ws.Range("E43") "  +1.18%  "
Set-TextValue $ws.Range("D44") "2.18"
Set-TextValue $ws.Range("E44") "  +4.46%  "
Set-TextValue $ws.Range("D45") "1.77"
Set-TextValue $ws.Range("E45") "  +0.42%  "
Set-TextValue $ws.Range("D46") "5.27"
Set-TextValue $ws.Range("E46") "  -0.45%  "
Set-TextValue $ws.Range("D47") "1.732.07"
Set-TextValue $ws.Range("E47") "  -0.18%  "
Set-TextValue $ws.Range("D48") "86.18"
Set-TextValue $ws.Range("E48") "  -2.01%  "
Set-TextValue $ws.Range("E49") "  +1.00%  "
Set-TextValue $ws.Range("D50") "0.0969"
Set-TextValue $ws.Range("E50") "  -0.06%  "
Set-TextValue $ws.Range("D51") "0.0497"
Set-TextValue $ws.Range("E51") "  -0.16%  "
